# ---------------------------------------------------------------------------
# Update metrics and price data files
#   - metrics_by_year: append FY2025 summary row
#   - signals_filtered: append new weekly/monthly signal rows through 2025,
#     add a couple of formula-driven return cells, freeze the header panes,
#     and make metrics_by_year the active sheet again.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)   # metrics_by_year
$sheet2 = $wb.Worksheets.Item(2)   # signals_filtered

# ---------------------------------------------------------------------------
# 1. metrics_by_year: new row for 2025
# ---------------------------------------------------------------------------
$sheet1.Range("A13").Value = 2025
$sheet1.Range("B13").Value = 15
$sheet1.Range("C13").Value = 53.33
$sheet1.Range("D13").Value = -0.603
$sheet1.Range("E13").Value = 89605.4141

# ---------------------------------------------------------------------------
# 2. signals_filtered: append new signal rows (61-72)
# ---------------------------------------------------------------------------
$newRows = @(
    @(61, 45712, "1wk", 3, 30,  94248.3515625,  2025, 89605.414099999995, -0.0492627979484721,  "no"),
    @(62, 45719, "1wk", 3, 20,  80601.0390625,  2025, 89605.414099999995,  0.11171537122515,     "no"),
    @(63, 45740, "1wk", 3, 30,  82334.5234375,  2025, 89605.414099999995,  0.08830913596675294,  "no"),
    @(64, 45747, "1wk", 3, 20,  78214.484375,   2025, 89605.414099999995,  0.1456370872802293,   "no"),
    @(65, 45824, "1wk", 3, 30,  100987.140625,  2025, 89605.414099999995, -0.1127047116302091,   "no"),
    @(66, 45894, "1wk", 3, 30,  108236.7109375, 2025, 89605.414099999995, -0.17213472872211,     "no"),
    @(67, 45962, "1mo", 3, 30,  90394.3125,     2025, 89605.414099999995, -0.008375889937765724, "no"),
    @(68, 45964, "1wk", 3, 30,  104719.640625,  2025, 89605.414099999995, -0.1443303899086504,   "no"),
    @(69, 45971, "1wk", 3, 20,  94177.078125,   2025, 89605.414099999995, -0.04854327776480905,  "no"),
    @(70, 45978, "1wk", 3, 10,  86805.0078125,  2025, 89605.414099999995,  0.03226088356617529,  "no"),
    @(71, 45985, "1wk", 3, 20,  90376.75,       2025, 89605.414099999995,  0.03226088356617529,  "no"),
    @(72, 45992, "1mo", 3, 20,  88430.1328125,  2025, 89605.414099999995,  0.01364972364747352,  "no")
)

foreach ($r in $newRows) {
    $row = $r[0]
    $sheet2.Cells.Item($row, 1).Value = $r[1]
    $sheet2.Cells.Item($row, 1).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
    $sheet2.Cells.Item($row, 2).Value = $r[2]
    $sheet2.Cells.Item($row, 3).Value = $r[3]
    $sheet2.Cells.Item($row, 4).Value = $r[4]
    $sheet2.Cells.Item($row, 5).Value = $r[5]
    $sheet2.Cells.Item($row, 6).Value = $r[6]
    $sheet2.Cells.Item($row, 7).Value = $r[7]
    $sheet2.Cells.Item($row, 8).Value = $r[8]
    $sheet2.Cells.Item($row, 9).Value = $r[9]
}

# Rows 73-75: formula-driven return column, with 74/75 sharing one formula
$sheet2.Range("A73").Value = 45999
$sheet2.Range("A73").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$sheet2.Range("B73").Value = "1wk"
$sheet2.Range("C73").Value = 3
$sheet2.Range("D73").Value = 20
$sheet2.Range("E73").Value = 90162.91
$sheet2.Range("F73").Value = 2025
$sheet2.Range("G73").Value = 89605.414099999995
$sheet2.Range("H73").Formula = "=(G73-E73)/E73"
$sheet2.Range("I73").Value = "no"

$sheet2.Range("A74").Value = 46006
$sheet2.Range("A74").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$sheet2.Range("B74").Value = "1wk"
$sheet2.Range("C74").Value = 3
$sheet2.Range("D74").Value = 20
$sheet2.Range("E74").Value = 88230.77
$sheet2.Range("F74").Value = 2025
$sheet2.Range("G74").Value = 89605.414099999995
$sheet2.Range("I74").Value = "no"

$sheet2.Range("A75").Value = 46013
$sheet2.Range("A75").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$sheet2.Range("B75").Value = "1wk"
$sheet2.Range("C75").Value = 3
$sheet2.Range("D75").Value = 20
$sheet2.Range("E75").Value = 88577.42
$sheet2.Range("F75").Value = 2025
$sheet2.Range("G75").Value = 89605.414099999995
$sheet2.Range("I75").Value = "no"

$sheet2.Range("H74:H75").Formula = "=(G74-E74)/E74"

# Stray underline-styled (empty) cell introduced alongside the new data
$sheet2.Range("L64").Font.Underline = $true

# ---------------------------------------------------------------------------
# 3. View state: freeze panes on signals_filtered, re-select cells, set
#    page orientation, and re-activate metrics_by_year as the shown tab.
# ---------------------------------------------------------------------------
$sheet2.Activate()
$sheet2.Range("J17").Select()
$excel.ActiveWindow.FreezePanes = $true
$sheet2.Range("H64").Select()
$sheet2.PageSetup.Orientation = 1

$sheet1.Activate()
$sheet1.Range("G14").Select()
